$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two extra pet rows (old rows 4 and 5) so only the header
# plus two data rows remain.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# Row 2: pet #0 -> Draco the Dino🦖 (owned by Riccardo, age recalculated)
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = "Draco"
$ws.Range("C2").Value = "Dino🦖"
$ws.Range("E2").Value = "0 y, 0 m, and 1 d"
$ws.Range("F2").Value = "Hasty"
$ws.Range("G2").Value = "Riccardo"
$ws.Range("I2").Value = 0

# Row 3: pet #1 -> Tobee the Bee🐝 (owned by Danni, age recalculated)
$ws.Range("B3").Value = "Tobee"
$ws.Range("C3").Value = "Bee🐝"
$ws.Range("E3").Value = "0 y, 0 m, and 1 d"
$ws.Range("F3").Value = "Hardy"
$ws.Range("G3").Value = "Danni"
$ws.Range("I3").Value = 0
